$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New player rows (16-19) appended to the batsman table
$newRows = @(
    @{ Row = 16; Player = "JP Inglis";     Mat = 2; Inns = 2; NO = 0; Runs = 53;  HS = 49; Ave = 26.5;  BF = 41;  SR = 129.26; C100 = 0; C50 = 0; Zeros = 0; Fours = 3;  Sixes = 2;  Type = "BAT"; Pos = 5 },
    @{ Row = 17; Player = "MP Breetzke";   Mat = 9; Inns = 7; NO = 0; Runs = 117; HS = 33; Ave = 16.71; BF = 113; SR = 103.53; C100 = 0; C50 = 0; Zeros = 0; Fours = 11; Sixes = 3;  Type = "BAT"; Pos = 1 },
    @{ Row = 18; Player = "RD Rickelton";  Mat = 8; Inns = 8; NO = 1; Runs = 336; HS = 89; Ave = 48;    BF = 188; SR = 178.72; C100 = 0; C50 = 3; Zeros = 0; Fours = 33; Sixes = 20; Type = "BAT"; Pos = 2 },
    @{ Row = 19; Player = "BJ Jacobs";     Mat = 3; Inns = 3; NO = 1; Runs = 26;  HS = 18; Ave = 13;    BF = 27;  SR = 96.29;  C100 = 0; C50 = 0; Zeros = 0; Fours = 2;  Sixes = 1;  Type = "BAT"; Pos = 5 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Player
    $ws.Cells.Item($row, 2).Value = $r.Mat
    $ws.Cells.Item($row, 3).Value = $r.Inns
    $ws.Cells.Item($row, 4).Value = $r.NO
    $ws.Cells.Item($row, 5).Value = $r.Runs
    $ws.Cells.Item($row, 6).Value = $r.HS
    $ws.Cells.Item($row, 7).Value = $r.Ave
    $ws.Cells.Item($row, 8).Value = $r.BF
    $ws.Cells.Item($row, 9).Value = $r.SR
    $ws.Cells.Item($row, 10).Value = $r.C100
    $ws.Cells.Item($row, 11).Value = $r.C50
    $ws.Cells.Item($row, 12).Value = $r.Zeros
    $ws.Cells.Item($row, 13).Value = $r.Fours
    $ws.Cells.Item($row, 14).Value = $r.Sixes
    $ws.Cells.Item($row, 15).Value = $r.Type
    $ws.Cells.Item($row, 16).Value = $r.Pos
}

# Formatting to mirror the source workbook's styling for the appended rows
$ws.Range("A16:A19").Font.Color = 4540992
$ws.Range("B16:N18").Font.Color = 4540992
$ws.Range("B19:M19").Font.Color = 2236962

$ws.Range("B19:M19").Interior.Color = 16777215

$ws.Range("A1:P19").Select()
$ws.Range("J11").Select()
